$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.868.28'
$ws.Range("E2").Value = '  -1.43%  '
$ws.Range("D3").Value = '1.638.71'
$ws.Range("E3").Value = '  -1.11%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9992'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.95%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '215.33'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.58%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5023'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.25%  '
$ws.Range("E7").Value = '  -0.55%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2568'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.25%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06373'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.29%  '
$ws.Range("E10").Value = '  -1.94%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07736'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.40%  '
$ws.Range("D12").Value = '1.641.20'
$ws.Range("E12").Value = '  -0.33%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.251'
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Value = '1.864.55'
$ws.Range("E14").Value = '  -0.97%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.5454'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.43%  '
$ws.Range("D16").Value = '0.0₅7890'
$ws.Range("E16").Value = '  -1.91%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '64.02'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.32%  '
$ws.Range("D18").Value = '25.893.31'
$ws.Range("E18").Value = '  -1.31%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.003'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.41%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '202.11'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -4.01%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.387'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.74%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.877'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.26%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.972'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.95%  '
$ws.Range("E24").Value = '  -0.34%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.861'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +3.21%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '140.72'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.87%  '
$ws.Range("E27").Value = '  -3.83%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.62'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.68%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.766'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -3.64%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.242'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.25%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.04964'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.90%  '
$ws.Range("E32").Value = '  -3.25%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.190'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.52%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.542'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.29%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.369'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.49%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.620'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -4.36%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.8896'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -3.46%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.5628'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.71%  '
$ws.Range("D39").Value = '1.141.85'
$ws.Range("E39").Value = '  -2.23%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01562'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.85%  '
$ws.Range("E41").Value = '  -0.35%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.667'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.90%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '99.93'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.47%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.8065'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.22%  '
$ws.Range("D45").Value = '1.776.38'
$ws.Range("E45").Value = '  -1.10%  '
$ws.Range("D46").Value = '0.0₈116'
$ws.Range("E46").Value = '  -0.29%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.4522'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.68%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.002'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.36%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '54.74'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.36%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.05046'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.73%  '
$ws.Range("E51").Value = '  -0.46%  '
